$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dialoguetest")
$ws.Activate()

# "main_talk" rows were renamed to "talk" (Type column)
$ws.Range("C4").Value = "talk"
$ws.Range("C7").Value = "talk"

# "chr1_image"/"chr2_image" commands renamed to "..._change" (Cmd column)
$ws.Range("H23").Value = "chr1_image_change"
$ws.Range("H24").Value = "chr2_image_change"

$ws.Range("H16").Select()
